$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.365.66', '  +1.27%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.011.72', '  +4.68%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.05%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '325.03', '  +1.36%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  +0.07%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5130', '  +1.49%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.4260', '  +5.39%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08718', '  +5.07%  '),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.137', '  +2.96%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.68', '  +2.51%  '),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.011.69', '  +5.01%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.584', '  +2.79%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.465', '  +3.07%  '),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.14%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '94.33', '  +2.16%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001115', '  +1.46%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06529', '  +0.40%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '18.90', '  +3.84%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.07%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.204', '  +4.33%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.413.76', '  +1.32%  '),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.82', '  +4.48%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.268', '  +3.31%  '),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.248.01', '  +5.30%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '22.43', '  +0.93%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '162.24', '  -0.19%  '),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.426', '  +5.09%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '131.09', '  +1.44%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.140', '  +0.49%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1052', '  +1.51%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.083', '  +2.05%  '),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.828', '  -0.13%  '),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.371', '  +14.57%  '),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02525', '  +2.90%  '),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06676', '  +3.97%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.460', '  +0.73%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '12.40', '  +8.83%  '),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.137', '  +4.62%  '),
    @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2194', '  +1.80%  '),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6655', '  +2.97%  '),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.240', '  +2.28%  '),
    @('Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.000', '  +0.10%  '),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '13.70', '  +3.01%  '),
    @('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6164', '  +1.82%  '),
    @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.192', '  -1.27%  '),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.668', '  +0.75%  '),
    @('EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.261', '  +4.39%  '),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '124.34', '  +1.72%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '80.63', '  +2.01%  '),
    @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06905', '  +1.56%  ')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
}
